# Remove the Aspose.Slides "evaluation" watermark text boxes that were
# present on every slide while the deck was produced under an
# unlicensed/evaluation build of Aspose.Slides for .NET. Now that the
# examples ship with a licensed build, these watermark shapes (each named
# "TextBox", containing "Evaluation only. / Created with Aspose.Slides
# for .NET 2.0 14.1.2.0 / Copyright 2004-2014 Aspose Pty Ltd.") are
# deleted from every slide in the deck.

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = $slide.Shapes.Count; $shi -ge 1; $shi--) {
        $shape = $slide.Shapes.Item($shi)

        if ($shape.Name -eq "TextBox") {
            $isWatermark = $false

            if ($shape.HasTextFrame) {
                $text = $shape.TextFrame.TextRange.Text
                if ($text -like "*Evaluation only*" -or $text -like "*Aspose*") {
                    $isWatermark = $true
                }
            }

            if ($isWatermark) {
                $shape.Delete()
            }
        }
    }
}
